$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "38.759.51"
$ws.Range("E2").Value = "  +2.90%  "
$ws.Range("D3").Value = "2.098.76"
$ws.Range("E3").Value = "  +3.09%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").Value = "'228.95"
$ws.Range("E5").Value = "  +0.87%  "
$ws.Range("D6").Value = "'0.613"
$ws.Range("E6").Value = "  +1.17%  "
$ws.Range("D7").Value = "'60.66"
$ws.Range("E7").Value = "  +1.67%  "
$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = "  -0.13%  "
$ws.Range("D9").Value = "'0.384"
$ws.Range("E9").Value = "  +2.37%  "
$ws.Range("D10").Value = "'0.0839"
$ws.Range("E10").Value = "  +0.91%  "
$ws.Range("D11").Value = "'0.104"
$ws.Range("E11").Value = "  -0.16%  "
$ws.Range("D12").Value = "2.396.17"
$ws.Range("E12").Value = "  +2.63%  "
$ws.Range("D13").Value = "'15.02"
$ws.Range("E13").Value = "  +4.42%  "
$ws.Range("D14").Value = "'22.33"
$ws.Range("E14").Value = "  +5.87%  "
$ws.Range("D15").Value = "'0.793"
$ws.Range("E15").Value = "  +2.84%  "
$ws.Range("D16").Value = "'5.47"
$ws.Range("E16").Value = "  -0.53%  "
$ws.Range("D17").Value = "2.078.20"
$ws.Range("E17").Value = "  +1.96%  "
$ws.Range("D18").Value = "38.610.16"
$ws.Range("E18").Value = "  +2.50%  "
$ws.Range("D19").Value = "'71.51"
$ws.Range("E19").Value = "  +3.24%  "
$ws.Range("D20").Value = "'6.06"
$ws.Range("E20").Value = "  +1.68%  "
$ws.Range("D21").Value = "0.0₃0836"
$ws.Range("E21").Value = "  +1.70%  "
$ws.Range("D22").Value = "'225.50"
$ws.Range("E22").Value = "  +0.75%  "
$ws.Range("E23").Value = "  -0.17%  "
$ws.Range("E24").Value = "  -0.21%  "
$ws.Range("D25").Value = "'2.34"
$ws.Range("E25").Value = "  +3.64%  "
$ws.Range("D26").Value = "'170.84"
$ws.Range("E26").Value = "  +1.79%  "
$ws.Range("D27").Value = "'9.51"
$ws.Range("E27").Value = "  +1.18%  "
$ws.Range("E28").Value = "  +6.81%  "
$ws.Range("D29").Value = "'19.24"
$ws.Range("E29").Value = "  +2.48%  "
$ws.Range("E30").Value = "  +8.94%  "
$ws.Range("E31").Value = "  +0.41%  "
$ws.Range("D32").Value = "'2.34"
$ws.Range("E32").Value = "  +4.38%  "
$ws.Range("D33").Value = "'4.77"
$ws.Range("E33").Value = "  +7.01%  "
$ws.Range("D34").Value = "'4.51"
$ws.Range("E34").Value = "  +3.21%  "
$ws.Range("D35").Value = "'0.0611"
$ws.Range("D36").Value = "'2.40"
$ws.Range("E36").Value = "  +2.26%  "
$ws.Range("D37").Value = "'6.38"
$ws.Range("E37").Value = "  -2.15%  "
$ws.Range("D38").Value = "'3.54"
$ws.Range("E38").Value = "  +4.20%  "
$ws.Range("D39").Value = "'0.999"
$ws.Range("E39").Value = "  -0.11%  "
$ws.Range("D40").Value = "'18.52"
$ws.Range("E40").Value = "  +3.30%  "
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").Value = "'101.61"
$ws.Range("E41").Value = "  +4.91%  "
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "1.544.50"
$ws.Range("E42").Value = "  +0.96%  "
$ws.Range("D43").Value = "'0.0222"
$ws.Range("E43").Value = "  +3.41%  "
$ws.Range("B44").Value = "Cronos"
$ws.Range("C44").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D44").Value = "'0.0927"
$ws.Range("E44").Value = "  +2.11%  "
$ws.Range("B45").Value = "HuobiToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D45").Value = "'2.81"
$ws.Range("E45").Value = "  -0.94%  "
$ws.Range("D46").Value = "'7.65"
$ws.Range("E46").Value = "  +8.69%  "
$ws.Range("D47").Value = "'4.11"
$ws.Range("E47").Value = "  -2.83%  "
$ws.Range("E48").Value = "  +0.95%  "
$ws.Range("D49").Value = "'1.04"
$ws.Range("E49").Value = "  +2.72%  "
$ws.Range("D50").Value = "'2.99"
$ws.Range("E50").Value = "  +2.04%  "
$ws.Range("D51").Value = "2.285.09"
$ws.Range("E51").Value = "  +2.68%  "
